# Bump the "Förändrad" (Changed) date in column C for every data row
# (rows 2-339) from 2023-09-08 (serial 45177) to 2023-09-09 (serial 45178).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSerial = 45177
$newSerial = 45178

$lastRow = 339

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
